$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Range("B10").Value = "-"
$ws.Range("C10").Value = "Ludoff-Maq. Term. FL"
$ws.Range("D10").Value = "[Emerson-Elet. Digit. Básica, -, -, -]"
$ws.Range("F10").Value = "-"

# Row 11
$ws.Range("B11").Value = "-"
$ws.Range("C11").Value = "Ludoff-Maq. Term. FL"
$ws.Range("D11").Value = "[Emerson-Elet. Digit. Básica, Vinicius-Ajustagem, -, -]"
$ws.Range("F11").Value = "[Claudinei-Des. Maq. Cad., -]"

# Row 12
$ws.Range("B12").Value = "-"
$ws.Range("D12").Value = "[Cleidson-Elet. Digit. Básica, Vinicius-Ajustagem, Carlos-Tornearia, Aline S. M.-Metalografia]"
$ws.Range("F12").Value = "[Claudinei-Des. Maq. Cad., -]"

# Row 14
$ws.Range("B14").Value = "[-, Emerson-Elet. Digit. Básica, -, -]"
$ws.Range("C14").Value = "-"
$ws.Range("D14").Value = "[Cleidson-Elet. Digit. Básica, Vinicius-Ajustagem, Carlos-Tornearia, Aline S. M.-Metalografia]"
$ws.Range("E14").Value = "[-, Nilton-Mec. Tec. Res. Mat]"
$ws.Range("F14").Value = "[Claudinei-Des. Maq. Cad., -]"

# Row 15
$ws.Range("B15").Value = "[-, Emerson-Elet. Digit. Básica, -, -]"
$ws.Range("C15").Value = "-"
$ws.Range("D15").Value = "[Cleidson-Elet. Digit. Básica, Vinicius-Ajustagem, Carlos-Tornearia, Aline S. M.-Metalografia]"
$ws.Range("E15").Value = "[Vinicius-Des. Maq. Cad., Nilton-Mec. Tec. Res. Mat]"
$ws.Range("F15").Value = "[Vinicius-Des. Maq. Cad., -]"

# Row 16
$ws.Range("B16").Value = "-"
$ws.Range("D16").Value = "[Cleidson-Elet. Digit. Básica, -, Carlos-Tornearia, Aline S. M.-Metalografia]"
$ws.Range("E16").Value = "[-, Nilton-Mec. Tec. Res. Mat]"
$ws.Range("F16").Value = "[Vinicius-Des. Maq. Cad., -]"
